$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "'33190"
$ws.Range("B3").Value = "'45967307"
$ws.Range("C3").Value = "Ricky"
$ws.Range("E3").Value = "'3147"
$ws.Range("A4").Value = "'37515"
$ws.Range("B4").Value = "'37069173"
$ws.Range("C4").Value = "詹toniii"
$ws.Range("E4").Value = "'3006"
$ws.Range("A5").Value = "'58987"
$ws.Range("A6").Value = "'83879"
$ws.Range("E6").Value = "'1588"
$ws.Range("A11").Value = "'21422"
$ws.Range("E11").Value = "'3471"
$ws.Range("A12").Value = "'4599"
$ws.Range("B12").Value = "'47131129"
$ws.Range("C12").Value = "NAM"
$ws.Range("E12").Value = "'4020"
$ws.Range("A13").Value = "'6823"
$ws.Range("B13").Value = "'53060417"
$ws.Range("C13").Value = "㊥老纳信耶稣"
$ws.Range("E13").Value = "'3975"
$ws.Range("A14").Value = "'8369"
$ws.Range("B14").Value = "'4756174"
$ws.Range("C14").Value = "純希です"
$ws.Range("E14").Value = "'3907"
$ws.Range("A15").Value = "'8480"
$ws.Range("B15").Value = "'49710892"
$ws.Range("C15").Value = "MMMMMMM"
$ws.Range("E15").Value = "'3902"
$ws.Range("A16").Value = "'13771"
$ws.Range("B16").Value = "'53520939"
$ws.Range("C16").Value = "㊥虎哥tiger"
$ws.Range("E16").Value = "'3699"
$ws.Range("A17").Value = "'24477"
$ws.Range("E17").Value = "'3394"
$ws.Range("A18").Value = "'47090"
$ws.Range("E18").Value = "'2714"
$ws.Range("A19").Value = "'55686"
$ws.Range("E19").Value = "'2553"
$ws.Range("A20").Value = "'8890"
$ws.Range("B20").Value = "'54698813"
$ws.Range("C20").Value = "閃亮唐老鴨"
$ws.Range("E20").Value = "'3885"
$ws.Range("A21").Value = "'9723"
$ws.Range("E21").Value = "'3851"
$ws.Range("A22").Value = "'10077"
$ws.Range("B22").Value = "'8057001"
$ws.Range("C22").Value = "㊥兵者诡道也"
$ws.Range("E22").Value = "'3836"
$ws.Range("A23").Value = "'10225"
$ws.Range("B23").Value = "'54085771"
$ws.Range("C23").Value = "㊥Matthieu"
$ws.Range("E23").Value = "'3830"
$ws.Range("A24").Value = "'11489"
$ws.Range("B24").Value = "'31495601"
$ws.Range("C24").Value = "陈晓军"
$ws.Range("E24").Value = "'3783"
$ws.Range("A25").Value = "'20695"
$ws.Range("E25").Value = "'3490"
$ws.Range("A26").Value = "'23683"
$ws.Range("E26").Value = "'3413"
$ws.Range("A27").Value = "'24342"
$ws.Range("E27").Value = "'3397"
$ws.Range("A28").Value = "'20780"
$ws.Range("B28").Value = "'3649043"
$ws.Range("C28").Value = "Dj6106"
$ws.Range("E28").Value = "'3488"
$ws.Range("A29").Value = "'30417"
$ws.Range("B29").Value = "'58408326"
$ws.Range("C29").Value = """Killer Bee"""
$ws.Range("E29").Value = "'3227"
$ws.Range("A30").Value = "'33036"
$ws.Range("E30").Value = "'3151"
$ws.Range("A31").Value = "'37216"
$ws.Range("E31").Value = "'3016"
$ws.Range("A32").Value = "'39127"
$ws.Range("E32").Value = "'2953"
$ws.Range("A33").Value = "'11119"
$ws.Range("E33").Value = "'3796"
$ws.Range("A34").Value = "'11686"
$ws.Range("E34").Value = "'3775"
$ws.Range("A35").Value = "'11939"
$ws.Range("B35").Value = "'56133764"
$ws.Range("C35").Value = "ustcarter"
$ws.Range("E35").Value = "'3765"
$ws.Range("A36").Value = "'13227"
$ws.Range("B36").Value = "'7852598"
$ws.Range("C36").Value = "seiji"
$ws.Range("E36").Value = "'3719"
$ws.Range("A37").Value = "'13768"
$ws.Range("B37").Value = "'55317038"
$ws.Range("C37").Value = "necman12345"
$ws.Range("E37").Value = "'3699"
$ws.Range("A38").Value = "'13776"
$ws.Range("B38").Value = "'49043337"
$ws.Range("C38").Value = "FanXiFang1976"
$ws.Range("E38").Value = "'3699"
$ws.Range("A39").Value = "'22076"
$ws.Range("B39").Value = "'11582001"
$ws.Range("C39").Value = "iMinatoX4"
$ws.Range("E39").Value = "'3455"
$ws.Range("A40").Value = "'22774"
$ws.Range("B40").Value = "'54778421"
$ws.Range("C40").Value = "Emma"
$ws.Range("E40").Value = "'3436"
$ws.Range("A41").Value = "'23105"
$ws.Range("B41").Value = "'47459684"
$ws.Range("C41").Value = "㊥阿闹切克闹"
$ws.Range("E41").Value = "'3428"
$ws.Range("A42").Value = "'24564"
$ws.Range("B42").Value = "'6809364"
$ws.Range("C42").Value = """Scorp IP"""
$ws.Range("E42").Value = "'3391"
$ws.Range("A43").Value = "'25077"
$ws.Range("B43").Value = "'56379103"
$ws.Range("C43").Value = "Globalking1001"
$ws.Range("E43").Value = "'3378"
$ws.Range("A44").Value = "'27648"
$ws.Range("B44").Value = "'56573048"
$ws.Range("C44").Value = "Xiaotian"
$ws.Range("E44").Value = "'3307"
$ws.Range("A45").Value = "'30819"
$ws.Range("B45").Value = "'20737010"
$ws.Range("C45").Value = "混着玩..."
$ws.Range("E45").Value = "'3216"
$ws.Range("A46").Value = "'32544"
$ws.Range("B46").Value = "'50837459"
$ws.Range("C46").Value = "NINE日"
$ws.Range("E46").Value = "'3166"
$ws.Range("A47").Value = "'32920"
$ws.Range("B47").Value = "'58203298"
$ws.Range("C47").Value = "权旨qua"
$ws.Range("E47").Value = "'3154"
$ws.Range("A48").Value = "'33325"
$ws.Range("B48").Value = "'57813281"
$ws.Range("C48").Value = "XAUEN"
$ws.Range("E48").Value = "'3143"
$ws.Range("A49").Value = "'41181"
$ws.Range("B49").Value = "'55634661"
$ws.Range("C49").Value = "Opalus"
$ws.Range("E49").Value = "'2885"
$ws.Range("A50").Value = "'41686"
$ws.Range("B50").Value = "'52997727"
$ws.Range("C50").Value = "larios"
$ws.Range("E50").Value = "'2870"
$ws.Range("A51").Value = "'41845"
$ws.Range("B51").Value = "'59020292"
$ws.Range("C51").Value = "Sharnoth"
$ws.Range("E51").Value = "'2864"
$ws.Range("A52").Value = "'43934"
$ws.Range("B52").Value = "'32316256"
$ws.Range("C52").Value = """秋の風 .."""
$ws.Range("E52").Value = "'2804"
$ws.Range("A53").Value = "'43927"
$ws.Range("E53").Value = "'2804"
$ws.Range("A54").Value = "'52985"
$ws.Range("E54").Value = "'2590"
$ws.Range("A55").Value = "'56722"
$ws.Range("B55").Value = "'31401481"
$ws.Range("C55").Value = "Player-31401481"
$ws.Range("E55").Value = "'2541"
$ws.Range("A56").Value = "'58565"
$ws.Range("B56").Value = "'37861953"
$ws.Range("C56").Value = """Durex ๑• . •๑"""
$ws.Range("E56").Value = "'2523"
$ws.Range("A57").Value = "'61726"
$ws.Range("A58").Value = "'62411"
$ws.Range("A60").Value = "'51199"
$ws.Range("E60").Value = "'2619"
$ws.Range("A64").Value = "'25162"
$ws.Range("E64").Value = "'3376"
$ws.Range("A65").Value = "'38889"
$ws.Range("E65").Value = "'2961"
$ws.Range("A66").Value = "'64286"
$ws.Range("E66").Value = "'2491"
$ws.Range("A67").Value = "'41555"
$ws.Range("E67").Value = "'2874"
$ws.Range("A68").Value = "'52591"
$ws.Range("A69").Value = "'55477"
$ws.Range("E69").Value = "'2555"
$ws.Range("A70").Value = "'57003"
$ws.Range("E70").Value = "'2538"
$ws.Range("A71").Value = "'65339"
$ws.Range("E71").Value = "'2483"
$ws.Range("A72").Value = "'90742"
$ws.Range("E72").Value = "'1499"
$ws.Range("A86").Value = "'52983"
$ws.Range("E86").Value = "'2590"
$ws.Range("A87").Value = "'37873"
$ws.Range("E87").Value = "'2995"
$ws.Range("A89").Value = "'61419"
$ws.Range("E89").Value = "'2501"
$ws.Range("A91").Value = "'90415"
